# Update the extinction-trial image filenames to use Windows-style
# backslash path separators instead of forward slashes (the workbook was
# moved from a Mac path to a Windows path: "C:\Users\yhuang\Desktop\...").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2:A11 all hold the shared string "Extinction/CS+1.BMP"
$ws.Range("A2:A11").Value = "Extinction\CS+1.BMP"
# A12:A21 all hold the shared string "Extinction/CS-1.BMP"
$ws.Range("A12:A21").Value = "Extinction\CS-1.BMP"

# Move the active selection from A2 down to A20:B21 (active cell A20),
# matching where the user was last working in the sheet.
$ws.Range("A20:B21").Select()
